$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 86 (OCTOVENT PLUS SYRUP 100 ML): balance/transactions ratios change,
#    and the selling price text doubles from 29.0000 to 58.0000.
# ---------------------------------------------------------------------------
$ws.Range("H86").Value = "2:0"
$ws.Range("Q86").Value = "2:0"

$fmtP86 = $ws.Range("P86").NumberFormat
$ws.Range("P86").NumberFormat = "@"
$ws.Range("P86").Value = "58.0000"
$ws.Range("P86").NumberFormat = $fmtP86

# ---------------------------------------------------------------------------
# 2) Insert a new product row above row 137 ("معجون حلاقه 55555") for the
#    new item "مخمريه بلوب", pushing the following rows (old 137-140, the
#    totals row and the footer row) down by one.
# ---------------------------------------------------------------------------
$ws.Rows("137").Insert()

# Pull the (now shifted-down) neighbour row's cell formatting onto the blank
# row so the new row matches the table's look (borders/fill/number formats).
$ws.Range("A138:Q138").Copy()
$ws.Range("A137:Q137").PasteSpecial(-4122)
$ws.Range("A1").Select()

$ws.Rows("137").RowHeight = 25.5

$ws.Range("A137").Value2 = 131
$ws.Range("C137").Value = "مخمريه بلوب"
$ws.Range("H137").Value = "0:0"

$fmtL137 = $ws.Range("L137").NumberFormat
$ws.Range("L137").NumberFormat = "@"
$ws.Range("L137").Value = "0"
$ws.Range("L137").NumberFormat = $fmtL137

$ws.Range("N137").Value = "60.00"

$fmtP137 = $ws.Range("P137").NumberFormat
$ws.Range("P137").NumberFormat = "@"
$ws.Range("P137").Value = "60.0000"
$ws.Range("P137").NumberFormat = $fmtP137

$ws.Range("Q137").Value = "1:0"

$ws.Range("A137:B137").Merge()
$ws.Range("C137:G137").Merge()
$ws.Range("H137:K137").Merge()
$ws.Range("L137:M137").Merge()
$ws.Range("N137:O137").Merge()

# ---------------------------------------------------------------------------
# 3) The totals row (now at row 142) must reflect the sum of the updated
#    selling-price column, including the newly inserted row.
# ---------------------------------------------------------------------------
$ws.Range("P142").Value2 = 8019.875
